# Insert a new weekly price row at row 84 (pushing existing rows 84..211 down to 85..212)
# and populate it with the new week's data for Orégano, Mercado Mayorista Lo Valledor de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 84; everything from 84 downward shifts to 85 downward.
$ws.Rows.Item(84).Insert()

# Fill in the new row 84 with the new data point.
$ws.Range("A84").Value = 6
$ws.Range("B84").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44721
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100112029
$ws.Range("G84").Value = "Orégano"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 49
$ws.Range("K84").Value = 12000
$ws.Range("L84").Value = 13000
$ws.Range("M84").Value = 12469
$ws.Range("N84").Value = "$/docena de atados"
$ws.Range("O84").Value = "Región Metropolitana"
$ws.Range("P84").Value = 4156
$ws.Range("Q84").Value = 3
$ws.Range("R84").Value = "Hortaliza"
